$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.8
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 1.08
$ws.Range("K3").Value = 8
$ws.Range("N3").Value = 2.3
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 1.7
$ws.Range("V3").Value = 9
$ws.Range("X3").Value = 17
$ws.Range("Y3").Value = 34
$ws.Range("AA3").Value = 7
$ws.Range("AB3").Value = 19
$ws.Range("AD3").Value = 501
$ws.Range("AE3").Value = 10
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 15

# Row 4
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 5.75
$ws.Range("U4").Value = 7
$ws.Range("W4").Value = 11
$ws.Range("AB4").Value = 19
$ws.Range("AD4").Value = 351
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 29
$ws.Range("AH4").Value = 67

# Row 5
$ws.Range("G5").Value = 2.45
$ws.Range("I5").Value = 3.3
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 9.5
$ws.Range("AF5").Value = 15

# Row 6
$ws.Range("J6").Value = 1.07
$ws.Range("K6").Value = 9
$ws.Range("L6").Value = 1.36
$ws.Range("M6").Value = 3.2
$ws.Range("N6").Value = 2.1
$ws.Range("O6").Value = 1.7

# Row 7
$ws.Range("G7").Value = 3.75
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.07
$ws.Range("J7").Value = 1.07
$ws.Range("K7").Value = 7.5
$ws.Range("L7").Value = 1.38
$ws.Range("M7").Value = 2.9
$ws.Range("N7").Value = 2.18
$ws.Range("Q7").Value = 2.57
$ws.Range("R7").Value = 1.88
$ws.Range("S7").Value = 1.83
$ws.Range("T7").Value = 10
$ws.Range("U7").Value = 21
$ws.Range("V7").Value = 13
$ws.Range("W7").Value = 60
$ws.Range("X7").Value = 37
$ws.Range("Y7").Value = 45
$ws.Range("Z7").Value = 7.5
$ws.Range("AA7").Value = 6.2
$ws.Range("AB7").Value = 15.5
$ws.Range("AC7").Value = 80
$ws.Range("AD7").Value = 700
$ws.Range("AE7").Value = 6.2
$ws.Range("AF7").Value = 9.25
$ws.Range("AG7").Value = 9

# Row 9
$ws.Range("G9").Value = 2.18
$ws.Range("H9").Value = 3.3
$ws.Range("J9").Value = 1.09
$ws.Range("K9").Value = 6.6
$ws.Range("L9").Value = 1.42
$ws.Range("M9").Value = 2.72
$ws.Range("N9").Value = 2.25
$ws.Range("O9").Value = 1.6
$ws.Range("P9").Value = 1.47
$ws.Range("Q9").Value = 2.55
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 1.75
$ws.Range("T9").Value = 6.3
$ws.Range("U9").Value = 10
$ws.Range("V9").Value = 10
$ws.Range("X9").Value = 22
$ws.Range("Y9").Value = 40
$ws.Range("Z9").Value = 6.6
$ws.Range("AA9").Value = 6.7
$ws.Range("AB9").Value = 19
$ws.Range("AC9").Value = 110
$ws.Range("AD9").Value = 1250
$ws.Range("AE9").Value = 8.25
$ws.Range("AF9").Value = 16.5
$ws.Range("AG9").Value = 12.5
$ws.Range("AI9").Value = 35
$ws.Range("AJ9").Value = 50

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 3.25
$ws.Range("P10").Value = 1.57
$ws.Range("Q10").Value = 2.25
$ws.Range("T10").Value = 6
$ws.Range("W10").Value = 21
$ws.Range("Z10").Value = 6.5
$ws.Range("AB10").Value = 19
$ws.Range("AG10").Value = 13
$ws.Range("AI10").Value = 34

# Row 11
$ws.Range("G11").Value = 3.9
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 1.93
$ws.Range("L11").Value = 1.28
$ws.Range("M11").Value = 3
$ws.Range("N11").Value = 1.91
$ws.Range("P11").Value = 1.42
$ws.Range("Q11").Value = 2.47
$ws.Range("R11").Value = 1.7
$ws.Range("S11").Value = 1.91
$ws.Range("T11").Value = 11.75
$ws.Range("U11").Value = 24
$ws.Range("W11").Value = 70
$ws.Range("X11").Value = 37
$ws.Range("Z11").Value = 8.75
$ws.Range("AA11").Value = 6
$ws.Range("AB11").Value = 13
$ws.Range("AC11").Value = 60
$ws.Range("AD11").Value = 450
$ws.Range("AE11").Value = 6.7
$ws.Range("AF11").Value = 9.25
$ws.Range("AG11").Value = 8.25
$ws.Range("AH11").Value = 18
$ws.Range("AI11").Value = 16
$ws.Range("AJ11").Value = 27

# Row 12
$ws.Range("G12").Value = 1.48
$ws.Range("H12").Value = 4.5
$ws.Range("I12").Value = 6
$ws.Range("N12").Value = 1.65
$ws.Range("O12").Value = 2.2
$ws.Range("P12").Value = 1.3
$ws.Range("Q12").Value = 3.4
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 1.91
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 7.5
$ws.Range("Y12").Value = 23
$ws.Range("Z12").Value = 15
$ws.Range("AB12").Value = 17
$ws.Range("AD12").Value = 201

# Row 27
$ws.Range("J27").Value = 1.1
$ws.Range("K27").Value = 7
$ws.Range("N27").Value = 2.4
$ws.Range("O27").Value = 1.53
$ws.Range("AD27").Value = 1250

# Row 40
$ws.Range("J40").Value = 1.07
$ws.Range("K40").Value = 9
